# Update column I ("10/03/2023" counts) and the dependent Delta_Offerto
# percentage in column J for the MOB p.2 chart prelude.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  I = 1225.9;  J = -99.86442613590015 },
    @{ Row = 3;  I = 605;     J = -19.83471074380165 },
    @{ Row = 4;  I = 13;      J = 61.53846153846154 },
    @{ Row = 5;  I = 235;     J = 91.48936170212767 },
    @{ Row = 6;  I = 105;     J = -35.23809523809523 },
    @{ Row = 7;  I = 212;     J = -66.98113207547169 },
    @{ Row = 10; I = 989;     J = -64.00404448938322 },
    @{ Row = 11; I = 400;     J = -1.000000000000001 },
    @{ Row = 12; I = 1116.7;  J = 1.101459657920656 },
    @{ Row = 13; I = 849;     J = -94.81743227326265 },
    @{ Row = 14; I = 1280;    J = -38.046875 },
    @{ Row = 15; I = 343;     J = -18.95043731778425 },
    @{ Row = 17; I = 208;     J = -54.32692307692308 },
    @{ Row = 19; I = 23;      J = 0 },
    @{ Row = 20; I = 126;     J = -49.20634920634921 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
